# Remove the trailing "Ver no Jupiter ..." and "© 2020 ..." paragraphs
# (plus the blank paragraph separating them from the requirements text)
# that used to be appended at the end of the page, right after the
# "LOT2013: Engenharia Bioquímica I (Requisito fraco)" requirement line.

$d = $word.ActiveDocument

$target1 = "Ver no Jupiter Salvar em pdf Salvar em docx"
$target2 = "Powered by Jekyll and Github pages"

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like ("*" + $target1 + "*")) {
        # The blank paragraph right before this one is part of the block
        # being removed too, so start deleting from there.
        $startPara = $d.Paragraphs.Item($i - 1)
    }
    if ($t -like ("*" + $target2 + "*")) {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $range.Delete()
}
